# "optimize quest_&_answ_randomizer.py + randomize answers"
# The generator script was re-run, producing a freshly randomized set of
# answers/questions for both quiz sheets, plus a couple of cosmetic view
# tweaks (column width, selection, active sheet) left over from the
# authoring session.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "exemple" (sheet1) - randomized answer order (question text
# itself is unchanged, only the "Reponse N" choices move / get re-tagged
# with their block suffix).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("C1").Value  = "Reponse 1-1"
$ws1.Range("C2").Value  = "Reponse 2-1"
$ws1.Range("C3").Value  = "Reponse 3-1"
$ws1.Range("C4").Value  = "Reponse 4-1"

$ws1.Range("C5").Value  = "Reponse 1-2"
$ws1.Range("C6").Value  = "Reponse 2-2"
$ws1.Range("C7").Value  = "Reponse 3-2"
$ws1.Range("C8").Value  = "Reponse 4-2"

$ws1.Range("C9").Value  = "Reponse 1-3"
$ws1.Range("C10").Value = "Reponse 2-3"
$ws1.Range("C11").Value = "Reponse 3-3"
$ws1.Range("C12").Value = "Reponse 4-3"

$ws1.Range("C13").Value = "Reponse 1-4"
$ws1.Range("C14").Value = "Reponse 2-4"
$ws1.Range("C15").Value = "Reponse 3-4"
$ws1.Range("C16").Value = "Reponse 4-4"

$ws1.Range("C17").Value = "Reponse 1-5"
$ws1.Range("C18").Value = "Reponse 2-5"
$ws1.Range("C19").Value = "Reponse 3-5"
$ws1.Range("C20").Value = "Reponse 4-5"

# widen the answer column and drop the stale selection left on D20
$ws1.Columns.Item(3).ColumnWidth = 25.140625

# ---------------------------------------------------------------------
# Sheet "exemple 2" (sheet2) - questions renumbered (Question 1-2 etc.
# became Question 12 etc.) and answers re-tagged with their block suffix.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("B1").Value  = "Question 12"
$ws2.Range("C1").Value  = "Reponse 12-1"
$ws2.Range("C2").Value  = "Reponse 22-1"
$ws2.Range("C3").Value  = "Reponse 32-1"
$ws2.Range("C4").Value  = "Reponse 42-1"

$ws2.Range("B5").Value  = "Question 22"
$ws2.Range("C5").Value  = "Reponse 12-2"
$ws2.Range("C6").Value  = "Reponse 22-2"
$ws2.Range("C7").Value  = "Reponse 32-2"
$ws2.Range("C8").Value  = "Reponse 42-2"

$ws2.Range("B9").Value  = "Question 33"
$ws2.Range("C9").Value  = "Reponse 12-3"
$ws2.Range("C10").Value = "Reponse 22-3"
$ws2.Range("C11").Value = "Reponse 32-3"
$ws2.Range("C12").Value = "Reponse 42-3"

$ws2.Range("B13").Value = "Question 44"
$ws2.Range("C13").Value = "Reponse 12-4"
$ws2.Range("C14").Value = "Reponse 22-4"
$ws2.Range("C15").Value = "Reponse 32-4"
$ws2.Range("C16").Value = "Reponse 42-4"

$ws2.Range("B17").Value = "Question 55"
$ws2.Range("C17").Value = "Reponse 12-5"
$ws2.Range("C18").Value = "Reponse 22-5"
$ws2.Range("C19").Value = "Reponse 32-5"
$ws2.Range("C20").Value = "Reponse 42-5"

# widen the answer column on sheet2 too
$ws2.Columns.Item(3).ColumnWidth = 17.42578125

# ---------------------------------------------------------------------
# View state: "exemple 2" ends up the active/selected tab with C23
# selected back on "exemple", and D18 selected on "exemple 2".
# ---------------------------------------------------------------------
$ws1.Range("C23").Select()
$ws2.Activate()
$ws2.Range("D18").Select()
